$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column F ("LabelR") holds HTML-ish tooltip strings for rows 2-7 that end
# with a trailing "<br>" which should be removed.
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().EndsWith("<br>")) {
        $cell.Value2 = $val.ToString().Substring(0, $val.ToString().Length - 4)
    }
}
